$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    7  = "https://www.lanacion.com.ar//politica/javier-milei-en-vivo-las-ultimas-medidas-del-gobierno-nid23032024/"
    8  = "https://www.lanacion.com.ar//politica/javier-milei-y-fatima-florez-fueron-al-teatro-y-cantaron-euforicos-canciones-de-sandro-nid24032024/"
    9  = "https://www.lanacion.com.ar//deportes/futbol/asi-quedo-el-cuadro-de-la-copa-argentina-2024-tras-el-triunfo-de-boca-nid24032024/"
    10 = "https://www.lanacion.com.ar//economia/el-consumidor-paga-fortunas-el-productor-cobra-migajas-nid24032024/"
    11 = "https://www.lanacion.com.ar//economia/las-rampas-de-acceso-y-las-cajas-transparentes-como-sigue-la-pelicula-de-la-inteligencia-artificial-nid24032024/"
    21 = "https://www.lanacion.com.ar//politica/javier-milei-y-fatima-florez-fueron-al-teatro-y-cantaron-euforicos-canciones-de-sandro-nid24032024/"
    83 = "https://www.lanacion.com.ar//politica/marcha-por-el-24-de-marzo-la-izquierda-el-kirchnerismo-y-grupos-de-derechos-humanos-se-movilizan-nid23032024/"
    84 = "https://www.lanacion.com.ar//politica/posse-extiende-su-influencia-a-la-afi-y-advierten-sobre-el-regreso-de-viejas-practicas-en-el-area-de-nid23032024/"
    85 = "https://www.lanacion.com.ar//economia/como-quedarian-los-haberes-jubilatorios-de-la-anses-en-abril-nid23032024/"
    86 = "https://www.lanacion.com.ar//seguridad/los-videos-claves-que-incriminan-a-la-empleada-domestica-del-ingeniero-asesinado-en-el-country-la-nid23032024/"
    87 = "https://www.lanacion.com.ar//el-mundo/la-historia-del-despiadado-comerciante-de-esclavos-que-dio-su-nombre-a-una-de-las-universidades-mas-nid23032024/"
}

foreach ($row in $values.Keys) {
    $ws.Range("A$row").Value = $values[$row]
}
